$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" strings (e.g. "1.005") look like plain numbers to
# Excel's normal cell-entry auto-detection, which would silently turn them
# into numeric values instead of the text Coinranking export expects.
# Forcing the NumberFormat to Text ("@") before the assignment keeps the
# literal text, and resetting the Style back to "Normal" afterwards avoids
# leaving a stray text-format style behind on the cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '23.677.36'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.627.92'
$ws.Range("E3").Value = '  -1.84%  '
Set-TextValue 'D4' '1.005'
$ws.Range("E4").Value = '  +0.34%  '
Set-TextValue 'D5' '307.64'
$ws.Range("E5").Value = '  -0.84%  '
Set-TextValue 'D6' '1.004'
$ws.Range("E6").Value = '  +0.21%  '
Set-TextValue 'D7' '0.3831'
$ws.Range("E7").Value = '  -1.88%  '
Set-TextValue 'D8' '0.3777'
$ws.Range("E8").Value = '  -2.50%  '
Set-TextValue 'D9' '50.47'
$ws.Range("E9").Value = '  -1.73%  '
Set-TextValue 'D10' '1.314'
$ws.Range("E10").Value = '  -4.06%  '
Set-TextValue 'D11' '1.005'
$ws.Range("E11").Value = '  +0.31%  '
Set-TextValue 'D12' '0.08330'
$ws.Range("E12").Value = '  -2.12%  '
Set-TextValue 'D13' '23.52'
$ws.Range("E13").Value = '  -2.12%  '
Set-TextValue 'D14' '6.900'
$ws.Range("E14").Value = '  -4.39%  '
Set-TextValue 'D15' '7.687'
$ws.Range("E15").Value = '  -4.14%  '
Set-TextValue 'D16' '0.00001286'
$ws.Range("E16").Value = '  -2.17%  '
$ws.Range("D17").Value = '1.618.71'
$ws.Range("E17").Value = '  -2.17%  '
Set-TextValue 'D18' '93.02'
$ws.Range("E18").Value = '  -1.80%  '
Set-TextValue 'D19' '0.06931'
$ws.Range("E19").Value = '  -0.92%  '
Set-TextValue 'D20' '19.24'
$ws.Range("E20").Value = '  -3.81%  '
Set-TextValue 'D21' '6.815'
$ws.Range("E21").Value = '  -2.65%  '
Set-TextValue 'D22' '1.003'
$ws.Range("E22").Value = '  +0.08%  '
Set-TextValue 'D23' '13.42'
$ws.Range("D24").Value = '23.678.91'
$ws.Range("E24").Value = '  -1.43%  '
Set-TextValue 'D25' '2.424'
$ws.Range("E25").Value = '  -2.60%  '
Set-TextValue 'D26' '2.839'
$ws.Range("E26").Value = '  -8.64%  '
Set-TextValue 'D27' '21.78'
$ws.Range("E27").Value = '  -2.25%  '
Set-TextValue 'D28' '152.06'
$ws.Range("E28").Value = '  -1.07%  '
Set-TextValue 'D29' '5.431'
$ws.Range("E29").Value = '  +2.13%  '
Set-TextValue 'D30' '136.06'
$ws.Range("E30").Value = '  -3.21%  '
Set-TextValue 'D31' '7.911'
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("D33").Value = '1.810.21'
$ws.Range("E33").Value = '  -1.38%  '
Set-TextValue 'D34' '0.9783'
$ws.Range("E34").Value = '  -6.95%  '
Set-TextValue 'D35' '0.07817'
$ws.Range("E35").Value = '  -4.42%  '
Set-TextValue 'D36' '0.02857'
$ws.Range("E36").Value = '  -4.80%  '
Set-TextValue 'D37' '6.517'
$ws.Range("E37").Value = '  -2.98%  '
Set-TextValue 'D38' '0.2626'
$ws.Range("E38").Value = '  -3.17%  '
Set-TextValue 'D39' '10.27'
$ws.Range("E39").Value = '  -7.92%  '
Set-TextValue 'D40' '0.09011'
$ws.Range("E40").Value = '  -1.57%  '
Set-TextValue 'D41' '0.7408'
$ws.Range("E41").Value = '  -2.41%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D42' '13.16'
$ws.Range("E42").Value = '  -3.47%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D43' '1.406'
$ws.Range("E43").Value = '  -1.20%  '
Set-TextValue 'D44' '16.37'
$ws.Range("E44").Value = '  -1.68%  '
Set-TextValue 'D45' '0.6813'
$ws.Range("E45").Value = '  -3.21%  '
Set-TextValue 'D46' '2.391'
$ws.Range("E46").Value = '  -4.56%  '
Set-TextValue 'D47' '4.049'
$ws.Range("E47").Value = '  -1.17%  '
Set-TextValue 'D48' '1.002'
$ws.Range("E48").Value = '  +0.08%  '
Set-TextValue 'D49' '0.08151'
$ws.Range("E49").Value = '  -2.23%  '
Set-TextValue 'D50' '133.19'
$ws.Range("E50").Value = '  -1.74%  '
Set-TextValue 'D51' '1.206'
$ws.Range("E51").Value = '  -2.63%  '
